$d = $word.ActiveDocument

# --- "Programa" section: split the numbered list run into 5 pieces with
#     manual line breaks (<w:br/>) between each numbered item. ---
$d.Content.Find.Execute("trabalho.2- A", $true, $false, $false, $false, $false, $true, 1, $false, "trabalho.^l2- A", 2) | Out-Null
$d.Content.Find.Execute("termodinâmico).3- A", $true, $false, $false, $false, $false, $true, 1, $false, "termodinâmico).^l3- A", 2) | Out-Null
$d.Content.Find.Execute("Gibbs-Helmholtz.4- Equilíbrio", $true, $false, $false, $false, $false, $true, 1, $false, "Gibbs-Helmholtz.^l4- Equilíbrio", 2) | Out-Null
$d.Content.Find.Execute("metal-óxido-O2(g).5- Equilíbrio", $true, $false, $false, $false, $false, $true, 1, $false, "metal-óxido-O2(g).^l5- Equilíbrio", 2) | Out-Null

# --- "Avaliação" section: split the grade-formula sentence from the
#     formula itself with a manual line break. ---
$d.Content.Find.Execute("expressão:NF=", $true, $false, $false, $false, $false, $true, 1, $false, "expressão:^lNF=", 2) | Out-Null

# --- "Bibliografia" section: split the concatenated reference list into
#     5 separate entries with manual line breaks between each. ---
$d.Content.Find.Execute("270-1.2) P. Atkins", $true, $false, $false, $false, $false, $true, 1, $false, "270-1.^l2) P. Atkins", 2) | Out-Null
$d.Content.Find.Execute("1600-9.3) S.Stolen", $true, $false, $false, $false, $false, $true, 1, $false, "1600-9.^l3) S.Stolen", 2) | Out-Null
$d.Content.Find.Execute("49230-6.4) R. DeHoff", $true, $false, $false, $false, $false, $true, 1, $false, "49230-6.^l4) R. DeHoff", 2) | Out-Null
$d.Content.Find.Execute("4065-9.5) Y.A. Chang", $true, $false, $false, $false, $false, $true, 1, $false, "4065-9.^l5) Y.A. Chang", 2) | Out-Null
